# The edit rotates the seven data records held in rows 4-10 of the
# "Artfynd" sheet: the two records that were in rows 9-10 move up to
# become rows 4-5, and the five records that were in rows 4-8 shift
# down by two rows to become rows 6-10. Every cell's own value travels
# with its record; nothing in rows 1-3 changes.
#
#   destination row -> source row
#   4  -> 9
#   5  -> 10
#   6  -> 4
#   7  -> 5
#   8  -> 6
#   9  -> 7
#   10 -> 8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 4
$lastRow = 10
$lastCol = 51   # column AY is the last used column

$sourceForDest = @{
    4  = 9
    5  = 10
    6  = 4
    7  = 5
    8  = 6
    9  = 7
    10 = 8
}

# ---- Snapshot every cell in the block (value + its runtime type) ----
# Read-before-write is required because several destination rows pull
# from rows that this same operation will also overwrite.
$snapVal = @{}
$snapType = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        $key = "$r,$c"
        $snapVal[$key] = $val
        if ($null -eq $val) {
            $snapType[$key] = "Empty"
        } else {
            $snapType[$key] = $val.GetType().Name
        }
    }
}

# ---- Write the rotated values back out ----
for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $sourceForDest[$destRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $key = "$srcRow,$c"
        $val = $snapVal[$key]
        $type = $snapType[$key]
        $cell = $ws.Cells.Item($destRow, $c)

        if ($type -eq "Empty") {
            # Source cell had no content -> make sure destination is blank too.
            $cell.Value2 = ""
        } elseif ($type -eq "String") {
            # Force text so Excel doesn't reinterpret numeric-looking or
            # date-looking text (e.g. "80", "2019-08-14") as a number/date.
            # A leading apostrophe is Excel's own "treat as text" marker;
            # restoring the original style afterwards drops the resulting
            # quote-prefix formatting flag so no stray style is left behind.
            $origStyle = $cell.Style
            $cell.Value2 = "'" + $val
            $cell.Style = $origStyle
        } else {
            # Numbers and booleans round-trip safely as-is.
            $cell.Value2 = $val
        }
    }
}

Write-Host "Row rotation complete"
